$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("KVO")

# Rename "index" column header to "i"
$ws.Range("A1").Value2 = "i"

# Shift index values from 1-based to 0-based (row 2 = A2:A503)
for ($r = 2; $r -le 503; $r++) {
    $ws.Cells.Item($r, 1).Value2 = $r - 2
}

# Update column A width (bestFit recalculation approximation)
$ws.Columns.Item(1).ColumnWidth = 3.1666666666666665

Write-Host "Done"
